$d = $word.ActiveDocument

# Helper: set the font size (in points) of an otherwise-empty paragraph
# (a paragraph that holds no runs, only paragraph-mark formatting).
# A direct "$para.Range.Font.Size = x" is a no-op on a zero-length range
# in this host, so we temporarily insert a placeholder character, apply
# the size to the (now non-empty) range, then delete the placeholder —
# leaving only the paragraph-mark run-properties changed, same as Word
# does when you place the cursor on a blank line and change the font
# size in the Ribbon/Font dialog.
function Set-EmptyParaFontSize($para, [double]$sizePt) {
    $r = $para.Range
    $r.InsertBefore("X")
    $r2 = $para.Range
    $r2.Font.Size = $sizePt
    $r2.Font.SizeBi = $sizePt
    $charRange = $d.Range($r2.Start, $r2.Start + 1)
    $charRange.Delete()
}

# Helper: turn an empty paragraph bold (same placeholder trick), used
# for the new bold spacer paragraph below the SiriusXM entry.
function Set-EmptyParaBold($para) {
    $r = $para.Range
    $r.InsertBefore("X")
    $r2 = $para.Range
    $r2.Font.Bold = $true
    $charRange = $d.Range($r2.Start, $r2.Start + 1)
    $charRange.Delete()
}

# --------------------------------------------------------------------
# 1) Shrink the thin spacer paragraphs from 8pt (sz 16) to 6pt (sz 12).
#    These are the blank paragraphs that sit directly above each of the
#    section headings that have a bottom border (Executive Summary,
#    Core Competencies, Work Experience, Education, Active
#    Certifications, Volunteering Experience).
# --------------------------------------------------------------------
Set-EmptyParaFontSize $d.Paragraphs(3) 6    # above "Executive Summary"
Set-EmptyParaFontSize $d.Paragraphs(7) 6    # above "Core Competencies"
Set-EmptyParaFontSize $d.Paragraphs(17) 6   # above "Work Experience"
Set-EmptyParaFontSize $d.Paragraphs(55) 6   # above "Education at the University of Texas at Austin"
Set-EmptyParaFontSize $d.Paragraphs(60) 6   # above "Active Certifications"
Set-EmptyParaFontSize $d.Paragraphs(64) 6   # above "Volunteering Experience"

# --------------------------------------------------------------------
# 2) Add a new, most-recent job entry (SiriusXM) at the top of the
#    "Work Experience" section, i.e. directly below the thin spacer
#    that follows the "Work Experience" heading (paragraph 19) and
#    above the existing "Amazon Web Services | Senior Solutions
#    Architect" entry.
# --------------------------------------------------------------------
$workSpacer = $d.Paragraphs(19)

# Insert three fresh, empty paragraphs right after the spacer; they
# pick up that spacer's (non-bold, sz 10) formatting, so each can be
# adjusted independently afterwards.
$workSpacer.Range.InsertParagraphAfter()
$workSpacer.Range.InsertParagraphAfter()
$workSpacer.Range.InsertParagraphAfter()

# Paragraph 20: bold title line "SiriusXM | Director | FinOps and Partnerships"
$titlePara = $d.Paragraphs(20)
$titlePara.Range.Text = "SiriusXM | Director | FinOps and Partnerships"
$titleRange = $d.Paragraphs(20).Range
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 11
$titleRange.Font.SizeBi = 11

# Paragraph 21: date/location line "June 2025 – present | Los Angeles, California"
$datePara = $d.Paragraphs(21)
$datePara.Range.Text = "June 2025 " + [char]0x2013 + " present | Los Angeles, California"
$dateRange = $d.Paragraphs(21).Range
$dateRange.Font.Size = 11
$dateRange.Font.SizeBi = 11

# Paragraph 22: thin bold spacer (sz 10) separating this entry from the
# next one, matching the pattern used between other job entries.
Set-EmptyParaBold $d.Paragraphs(22)

# --------------------------------------------------------------------
# 3) Mark the "Terraform Certified Associate" badge image run as
#    NoProof (<w:noProof/>), matching the other certification badges.
# --------------------------------------------------------------------
foreach ($shape in $d.InlineShapes) {
    if ($shape.AlternativeText -eq "Terraform Certified Associate") {
        $shape.Range.NoProofing = $true
    }
}

Write-Output "done"
